$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data values (MV -datos- refresh) matching the diff.
$updates = @{
    "D23" = 0.1
    "C26" = 80.8
    "D29" = 1.2
    "D32" = 2.4
    "D33" = 1.2
    "C41" = 97.2
    "D49" = 1.1
    "C50" = 102.9
    "D53" = 0.9
    "C56" = 105.8
    "D57" = 0.2
    "C59" = 106.4
    "D61" = 0.8
    "D62" = 1.2
    "C63" = 111.5
    "D63" = 1
    "C64" = 111
    "D64" = -0.5
    "C65" = 112.1
    "D65" = 1
    "C66" = 111.9
    "C67" = 113.4
    "D67" = 1.3
    "C68" = 114
    "D68" = 0.5
    "C69" = 109.5
    "D69" = -3.9
    "D70" = 2.1
    "C71" = 97.90000000000001
    "D71" = -12.4
    "D72" = 5
    "C73" = 109.3
    "D73" = 6.3
    "D74" = 3.6
    "C75" = 114.9
    "D75" = 1.4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
